# Applies an update to the "Artfynd" sheet: rows 3-9 are re-ordered/updated
# so that each row's identifying data (Id, Taxonsorteringsordning, Rödlistade,
# TaxonId, Artnamn, Vetenskapligt namn, Auktor, Ost, Nord and the optional
# "Publik kommentar" note) reflects a new arrangement, while the remaining
# columns (validation status, locality, dates, observer, etc.) stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each row, keyed by row number.
# Columns: A, B, D, E, F, G, H, Q, R, AC (AC optional - $null means "no comment")
$rowsData = @{
    3 = @{ A = 111639168; B = 89686; D = "NT"; E = 658;    F = "Rosenticka"; G = "Rhodofomes roseus"; H = "(Alb. & Schwein.) Kotl. & Pouzar"; Q = 548104.1391889038; R = 6926477.987023209; AC = $null }
    4 = @{ A = 111639167; B = 96348; D = "VU"; E = 220787; F = "Knärot";     G = "Goodyera repens";    H = "(L.) R. Br.";                      Q = 547814.5103353403; R = 6926124.461383951; AC = "1 planta" }
    5 = @{ A = 111639173; B = 96348; D = "VU"; E = 220787; F = "Knärot";     G = "Goodyera repens";    H = "(L.) R. Br.";                      Q = 547838.0352795018; R = 6926228.915831603; AC = "ca 15 plantor" }
    6 = @{ A = 111639172; B = 96348; D = "VU"; E = 220787; F = "Knärot";     G = "Goodyera repens";    H = "(L.) R. Br.";                      Q = 548221.3480213688; R = 6926511.607424877; AC = $null }
    7 = @{ A = 111639175; B = 89686; D = "NT"; E = 658;    F = "Rosenticka"; G = "Rhodofomes roseus";  H = "(Alb. & Schwein.) Kotl. & Pouzar"; Q = 547828.4099300706; R = 6926124.660841302; AC = $null }
    8 = @{ A = 111639170; B = 96348; D = "VU"; E = 220787; F = "Knärot";     G = "Goodyera repens";    H = "(L.) R. Br.";                      Q = 548231.4260436196; R = 6926519.619127685; AC = "ca 15 plantor" }
    9 = @{ A = 111639174; B = 96348; D = "VU"; E = 220787; F = "Knärot";     G = "Goodyera repens";    H = "(L.) R. Br.";                      Q = 547803.9854679118; R = 6926147.447742103; AC = "ca 6 plantor" }
}

foreach ($r in $rowsData.Keys) {
    $data = $rowsData[$r]

    $ws.Range("A$r").Value = $data.A
    $ws.Range("B$r").Value = $data.B
    $ws.Range("D$r").Value = $data.D
    $ws.Range("E$r").Value = $data.E
    $ws.Range("F$r").Value = $data.F
    $ws.Range("G$r").Value = $data.G
    $ws.Range("H$r").Value = $data.H
    $ws.Range("Q$r").Value = $data.Q
    $ws.Range("R$r").Value = $data.R

    if ($null -eq $data.AC) {
        $ws.Range("AC$r").ClearContents()
    } else {
        $ws.Range("AC$r").Value = $data.AC
    }
}
